$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F2: last activity timestamp
$ws.Range("F2").Value = "2024-10-18 14:00:55"

# G2: full name
$ws.Range("G2").Value = "Ntcn nt dsf"

# H2: phone number - force text so the leading "+" and digits are preserved
# as a string instead of being auto-converted to a number, then reset the
# cell style back to Normal so no extra number-format style is introduced.
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "+7965214251"
$ws.Range("H2").Style = "Normal"

# I2: email
$ws.Range("I2").Value = "asdaS@dasd.conm"

# K2: dog breed
$ws.Range("K2").Value = "sdf sdf"

# L2: dog name
$ws.Range("L2").Value = "sdgsfdgsdfg"

# M2: birth date - force text so it stays as the literal string "05.03.2024"
# instead of being parsed into a date serial number.
$ws.Range("M2").NumberFormat = "@"
$ws.Range("M2").Value = "05.03.2024"
$ws.Range("M2").Style = "Normal"

# N2: age
$ws.Range("N2").Value = "7 месяца(ев)"
